$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Java")
$ws.Range("B6").Value = 179
